# Weekly fruit/hortaliza update: insert a new weekly price record for
# "Pepino ensalada" (Vega Central Mapocho de Santiago) as row 235,
# shifting the existing rows 235-241 down to 236-242.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data down by inserting a blank row at position 235.
$ws.Rows.Item(235).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(235, 1).Value  = 9
$ws.Cells.Item(235, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(235, 3).Value  = "Metropolitana"
$ws.Cells.Item(235, 4).Value  = 44615
$ws.Cells.Item(235, 5).Value  = 13
$ws.Cells.Item(235, 6).Value  = 100112043
$ws.Cells.Item(235, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(235, 8).Value  = "Sin especificar"
$ws.Cells.Item(235, 9).Value  = "Primera"
$ws.Cells.Item(235, 10).Value = 61
$ws.Cells.Item(235, 11).Value = 14000
$ws.Cells.Item(235, 12).Value = 15000
$ws.Cells.Item(235, 13).Value = 14508
$ws.Cells.Item(235, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(235, 15).Value = "Región Metropolitana"
$ws.Cells.Item(235, 16).Value = 242
$ws.Cells.Item(235, 17).Value = 60
$ws.Cells.Item(235, 18).Value = "Hortaliza"

# Match the date-time style used by the other rows in column D.
$ws.Cells.Item(235, 4).NumberFormat = $ws.Cells.Item(236, 4).NumberFormat
